$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) and Volume(1h) (column E) cells are stored as text in
# this sheet (e.g. "27.265.83", "217.40", "  -0.33%  "). For the handful of
# Price values that happen to look like a plain number (e.g. "217.40"),
# a leading apostrophe forces Excel to keep them as text (quotePrefix),
# exactly like typing '217.40 into a General-formatted cell, instead of
# letting them be auto-converted to a numeric value.
$ws.Range("D2").Value = '27.265.83'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '1.647.92'
$ws.Range("E3").Value = '  -0.83%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = '''217.40'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").Value = '''0.0630'
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("D10").Value = '''20.05'
$ws.Range("E10").Value = '  +0.41%  '
$ws.Range("D11").Value = '''0.0844'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '1.877.54'
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").Value = '1.668.62'
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("E14").Value = '  -2.02%  '
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '''67.67'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").Value = '27.225.35'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("D18").Value = '0.0₃0741'
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("D19").Value = '''219.03'
$ws.Range("E19").Value = '  -1.44%  '
$ws.Range("E20").Value = '  -0.21%  '
$ws.Range("D21").Value = '''6.86'
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '''2.51'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  -0.99%  '
$ws.Range("D25").Value = '''147.47'
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  +1.67%  '
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("E28").Value = '  -0.93%  '
$ws.Range("D29").Value = '''15.81'
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").Value = '''0.0508'
$ws.Range("E30").Value = '  -1.40%  '
$ws.Range("E31").Value = '  -1.59%  '
$ws.Range("D32").Value = '''3.39'
$ws.Range("E32").Value = '  -0.13%  '
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("E34").Value = '  +1.31%  '
$ws.Range("D35").Value = '1.265.13'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("D38").Value = '''0.544'
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("E40").Value = '  -0.18%  '
$ws.Range("D41").Value = '''0.809'
$ws.Range("E41").Value = '  -0.73%  '
$ws.Range("D42").Value = '''2.23'
$ws.Range("E42").Value = '  +4.46%  '
$ws.Range("D43").Value = '''5.43'
$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("D44").Value = '1.787.92'
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("D45").Value = '''62.60'
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("D46").Value = '''91.99'
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  +13.76%  '
$ws.Range("E49").Value = '  -1.43%  '
$ws.Range("D50").Value = '''7.68'
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("E51").Value = '  -1.00%  '
